# "analise de algoritmo ordenação"
# Corrects the CASO MEDIO (average case) comparison-count table for the
# bubble / insertion / selection sort benchmarks on sheet "Plan1", and
# moves the active selection to where the author left off editing.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Plan1")

# --- Row 16 ("bubble"): fix the stuck O(n^2) figures for n = 5000/10000/50000
$ws.Range("D16").Value = 12497500
$ws.Range("E16").Value = 49995000
$ws.Range("F16").Value = 1249975000

# --- Row 17 ("insertion"): refreshed empirical comparison counts
$ws.Range("B17").Value = 62615
$ws.Range("C17").Value = 568079
$ws.Range("D17").Value = 6196419
$ws.Range("E17").Value = 24579512
$ws.Range("F17").Value = 626863600

# --- Row 18 ("selection"): same O(n^2) fix as row 16
$ws.Range("D18").Value = 12497500
$ws.Range("E18").Value = 49995000
$ws.Range("F18").Value = 1249975000

# --- Leave the cursor where the author finished up
$ws.Range("H15").Select()
